$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 new rows right after the header row (new rows 2-6),
#    shifting the existing data rows (old 2-21) down to rows 7-26.
$insertRange = $ws.Range("A2:C6")
$insertRange.Insert()
$insertRange.ClearFormats()

$topData = @(
    @(0.299476683139801, 0.1832595765590667, 0.4633412957191467),
    @(-0.08216137439012521, -0.3101668357849121, 0.4265366494655609),
    @(-0.1818851232528686, 0.0901026204228401, 0.0786489024758338),
    @(0.0897971913218498, 0.1569923609495163, 0.2553416788578033),
    @(0.1241583600640297, 0.194713294506073, 0.1459967941045761)
)

for ($i = 0; $i -lt $topData.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $topData[$i][0]
    $ws.Cells.Item($row, 2).Value = $topData[$i][1]
    $ws.Cells.Item($row, 3).Value = $topData[$i][2]
}

# 2) Append 5 new rows at the bottom (rows 27-31).
$bottomData = @(
    @(0.0612392425537109, 0.2464841306209564, 0.1278235465288162),
    @(0.0074830991216003, -0.0374154970049858, -0.0755945742130279),
    @(0.0713185146450996, 0.0108428578823804, 0.0226020142436027),
    @(0.0126754539087414, 0.0655152946710586, 0.0435241498053073),
    @(0.0337503030896186, -0.0120645882561802, -0.0161879286170005)
)

for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $row = 27 + $i
    $ws.Cells.Item($row, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($row, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($row, 3).Value = $bottomData[$i][2]
}
